$d = $word.ActiveDocument

# --- Special case: "Human Research Ethics Committee" appears twice in the
# document. Only the standalone occurrence inside the contacts table
# (the second one) should be translated; the first occurrence is embedded
# inside a longer English sentence that stays untouched. We locate the
# second occurrence via Find (without replacing) and then set its Range
# text directly.
$rng = $d.Content
$foundFirst = $rng.Find.Execute("Human Research Ethics Committee", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $foundFirst) { Write-Host "FAILED to find first Human Research Ethics Committee occurrence" }
$rng.Start = $rng.End
$foundSecond = $rng.Find.Execute("Human Research Ethics Committee", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $foundSecond) { Write-Host "FAILED to find second Human Research Ethics Committee occurrence" }
$rng.Text = "Ikomiti yeeNdlela zokuziphatha zoPhando loLuntu"

# --- Remaining unique text replacements ---

$result0 = $d.Content.Find.Execute("Who do I contact if I have questions or concerns?", $true, $false, $false, $false, $false, $true, 1, $false, "Ngubani endinokuqhagamshelana naye ukuba ndinemibuzo okanye iinkxalabo?", 2)
if (-not $result0) { Write-Host "FAILED replace #0: Who do I contact if I have questions or " }
$result1 = $d.Content.Find.Execute("If you have any questions or concerns about your rights as a study participant, you can contact the study team at ", $true, $false, $false, $false, $false, $true, 1, $false, "Ukuba unayo nayiphi na imibuzo okanye iinkxalabo malunga namalungelo akho njengomthathi-nxaxheba kuphononongo, ungaqhagamshelana neqela lophononongo ku-", 2)
if (-not $result1) { Write-Host "FAILED replace #1: If you have any questions or concerns ab" }
$result2 = $d.Content.Find.Execute(" or on WhatsApp at +27 XX XXX XXXX (messages only).", $true, $false, $false, $false, $false, $true, 1, $false, " okanye ku-WhatsApp ku +27 XX XXX XXXX (imiyalezo kuphela).", 2)
if (-not $result2) { Write-Host "FAILED replace #2:  or on WhatsApp at +27 XX XXX XXXX (mess" }
$result3 = $d.Content.Find.Execute("If you have more questions or concerns about your rights, you can contact one of the ethics committees listed: ", $true, $false, $false, $false, $false, $true, 1, $false, "Ukuba uneminye imibuzo okanye iinkxalabo malunga namalungelo akho, ungaqhagamshelana nenye yee komiti yokuziphatha edwelisiweyo: ", 2)
if (-not $result3) { Write-Host "FAILED replace #3: If you have more questions or concerns a" }
$result4 = $d.Content.Find.Execute("Name", $true, $false, $false, $false, $false, $true, 1, $false, "Igama", 2)
if (-not $result4) { Write-Host "FAILED replace #4: Name" }
$result5 = $d.Content.Find.Execute("Telephone", $true, $false, $false, $false, $false, $true, 1, $false, "Inombolo yomnxeba", 2)
if (-not $result5) { Write-Host "FAILED replace #5: Telephone" }
$result6 = $d.Content.Find.Execute("Email", $true, $false, $false, $false, $false, $true, 1, $false, "I-imeyile", 2)
if (-not $result6) { Write-Host "FAILED replace #6: Email" }
$result7 = $d.Content.Find.Execute("University of Cape Town Centre for Social Science Research ", $true, $false, $false, $false, $false, $true, 1, $false, "IDyunivesithi yaseKapa iZiko lezoPhando lweNzululwazi yeZentlalo ", 2)
if (-not $result7) { Write-Host "FAILED replace #7: University of Cape Town Centre for Socia" }
$result9 = $d.Content.Find.Execute("Informed Consent to Take Part in the Study", $true, $false, $false, $false, $false, $true, 1, $false, "Imvume eChaziweyo yokuThatha Inxaxheba kuPhononongo", 2)
if (-not $result9) { Write-Host "FAILED replace #9: Informed Consent to Take Part in the Stu" }
$result10 = $d.Content.Find.Execute("Please read these statements carefully: ", $true, $false, $false, $false, $false, $true, 1, $false, "Nceda ufunde ezi ngxelo ngononophelo:", 2)
if (-not $result10) { Write-Host "FAILED replace #10: Please read these statements carefully: " }
$result11 = $d.Content.Find.Execute("I have read the information above and know what is expected of my child.", $true, $false, $false, $false, $false, $true, 1, $false, "Ndilufundile olu lwazi lungentla kwaye ndiyayazi into ekufuneka yenziwe ngumntwana wam.", 2)
if (-not $result11) { Write-Host "FAILED replace #11: I have read the information above and kn" }
$result12 = $d.Content.Find.Execute("I understand as my child’s guardian that I am giving consent for them to participate. ", $true, $false, $false, $false, $false, $true, 1, $false, "Ndiyaqonda njengomgcini womntwana wam ukuba ndinika imvume yokuba athathe inxaxheba.", 2)
if (-not $result12) { Write-Host "FAILED replace #12: I understand as my child’s guardian that" }
$result13 = $d.Content.Find.Execute("I understand that even though I have given consent that my child will still be able to choose freely if they want to be interviewed. ", $true, $false, $false, $false, $false, $true, 1, $false, "Ndiyayiqonda ukuba nangona ndiyinikezile imvume umntwana wam usezokwazi ukukhetha ngokukhululekileyo ukuba uyalufuna udliwano-ndlebe.", 2)
if (-not $result13) { Write-Host "FAILED replace #13: I understand that even though I have giv" }
$result14 = $d.Content.Find.Execute("I understand that they can say no to being interviewed without any consequence. ", $true, $false, $false, $false, $false, $true, 1, $false, "Ndiyayiqonda ukuba banokuthi hayi kudliwano-ndlebe ngaphandle kweziphumo.", 2)
if (-not $result14) { Write-Host "FAILED replace #14: I understand that they can say no to bei" }
$result15 = $d.Content.Find.Execute("I had time to think about the information and have asked any questions I might have on either the email or message only WhatsApp number provided. I got satisfying answers if I did ask questions.", $true, $false, $false, $false, $false, $true, 1, $false, "Ndibe nexesha lokucinga malunga nolwazi kwaye ndibuze nayiphi na imibuzo endinokuba nayo kwi-imeyile okanye umyalezo kuphela kwinombolo ka-WhatsApp enikeziweyo. Ndifumene iimpendulo ezanelisayo xa bendibuza imibuzo.", 2)
if (-not $result15) { Write-Host "FAILED replace #15: I had time to think about the informatio" }
$result16 = $d.Content.Find.Execute("I know who can see my child’s information, how it will be kept safe, and what happens to it after the study.", $true, $false, $false, $false, $false, $true, 1, $false, "Ndiyayazi ukuba ngubani onokubona ulwazi lomntwana wam, ukuba luya kugcinwa njani lukhuselekile, kwaye kwenzeka ntoni kulo emva kophononongo.", 2)
if (-not $result16) { Write-Host "FAILED replace #16: I know who can see my child’s informatio" }
$result17 = $d.Content.Find.Execute("I understand that I will not be notified of my child’s answers.", $true, $false, $false, $false, $false, $true, 1, $false, "Ndiyayiqonda ukuba andizukwaziswa ngeempendulo zomntwana wam.", 2)
if (-not $result17) { Write-Host "FAILED replace #17: I understand that I will not be notified" }
$result18 = $d.Content.Find.Execute("I know that if the researchers pick up any safety concerns that they will let my child know that they will have to share the information before sharing it with me.", $true, $false, $false, $false, $false, $true, 1, $false, "Ndiyayazi ukuba ukuba abaphandi baye bafumanisa naziphi na iinkxalabo zokhuseleko bazakuxelela umntwana wam ayazi ukuba kuzakufuneka babelane ngolwazi ngaphambi kokuba babelane nam ngalo.", 2)
if (-not $result18) { Write-Host "FAILED replace #18: I know that if the researchers pick up a" }
$result19 = $d.Content.Find.Execute("I know I and my child won’t be named in any papers or reports from this study.", $true, $false, $false, $false, $false, $true, 1, $false, "Ndiyayazi ukuba mna nomntwana wam asisayi kuchazwa kuwo nawaphi na amaphepha okanye iingxelo zolu phononongo.", 2)
if (-not $result19) { Write-Host "FAILED replace #19: I know I and my child won’t be named in " }
$result20 = $d.Content.Find.Execute("I know who to contact if I have a problem with the study.", $true, $false, $false, $false, $false, $true, 1, $false, "Ndiyayazi ukuba ndiqhagamshelane nabani ukuba ndinengxaki ngophononongo.", 2)
if (-not $result20) { Write-Host "FAILED replace #20: I know who to contact if I have a proble" }
$result21 = $d.Content.Find.Execute("You can contact me again if more information is needed from me.", $true, $false, $false, $false, $false, $true, 1, $false, "Ungaphinda uqhagamshelane nam ukuba ulwazi oluninzi luyafuneka kum.", 2)
if (-not $result21) { Write-Host "FAILED replace #21: You can contact me again if more informa" }
$result22 = $d.Content.Find.Execute("You can keep my contact information safe so you can tell me about the results of the study.", $true, $false, $false, $false, $false, $true, 1, $false, "Ungazigcina iinkcukacha zam zoqhagamshelwano zikhuselekile ukuze undixelele ngeziphumo zophononongo.", 2)
if (-not $result22) { Write-Host "FAILED replace #22: You can keep my contact information safe" }
# Item #23 mixes a curly quote pair (around "Ewe"/"Yes") with a straight
# quote pair (around "Hayi"/"No"). Find.Execute's replacement text path
# runs through Word's smart-quote autocorrect and would turn the straight
# quotes into curly ones, so we locate the range first and then assign
# Range.Text directly, which performs a literal (non-autocorrected) swap.
$rng23 = $d.Content
$result23 = $rng23.Find.Execute('If you have read and understand the above document, agree with the messages and give consent to participate in the study, select “Yes” in WhatsApp. Select "No" in WhatsApp if you do not want to participate.', $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $result23) { Write-Host "FAILED replace #23: If you have read and understand the abov" }
$rng23.Text = 'Ukuba ulufundile kwaye waluqonda olu xwebhu lungasentla, uyavumelana nemiyalezo kwaye unike imvume yokuthatha inxaxheba kuphononongo, khetha u-“Ewe” kuWhatsApp. Khetha "Hayi" kuWhatsApp ukuba awufuni ukuthatha inxaxheba.'

Write-Host "All replacements attempted."
